$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# Row 6 previously held "Investor 5"; change it to "Investor 4" so it
# matches row 5 and the now-unused "Investor 5" shared string drops out.
$ws.Range("A6").Value = "Investor 4"

# Move the active selection to A7 (single cell) instead of A2:A6.
$ws.Range("A7").Select()
